$d = $word.ActiveDocument
$t = $d.Tables(1)

$values = @(
    @("48+7=", "92-3=", "4+67=", "9+58=", "57-8="),
    @("84-19=", "20-13=", "37+26=", "7+4=", "80-68="),
    @("32-17=", "62-48=", "83-18=", "32-28=", "96-38="),
    @("95-68=", "91-65=", "30-13=", "59+3=", "69+23="),
    @("87-18=", "2+39=", "86-29=", "41-24=", "59+33="),
    @("65-57=", "85-56=", "95-88=", "84-79=", "57+35="),
    @("52-46=", "37+29=", "42-14=", "81-28=", "77-28="),
    @("80-65=", "47-38=", "48-29=", "27+57=", "46+25="),
    @("60-35=", "60-46=", "86-38=", "93-57=", "80-5="),
    @("83-56=", "63-29=", "19+54=", "85-78=", "37+19="),
    @("70-15=", "6+87=", "27+15=", "49+36=", "58+4="),
    @("20-7=", "57+9=", "36+57=", "85-18=", "45-39="),
    @("48+48=", "92-74=", "18-9=", "28-19=", "27+14="),
    @("58+15=", "79+13=", "39+37=", "7+24=", "23-9="),
    @("82-77=", "31-18=", "47+45=", "33+39=", "16+46="),
    @("19+26=", "67-8=", "29+38=", "46+9=", "72-38="),
    @("69+8=", "64+9=", "70-17=", "34+57=", "69+28="),
    @("43-24=", "52-5=", "68+7=", "95-16=", "82-53="),
    @("67-29=", "48+34=", "48+14=", "72-53=", "60-49="),
    @("27+8=", "72-24=", "5+7=", "67+7=", "8+84=")
)

for ($r = 1; $r -le $values.Count; $r++) {
    $row = $values[$r - 1]
    for ($c = 1; $c -le $row.Count; $c++) {
        $t.Cell($r, $c).Range.Text = $row[$c - 1]
    }
}

Write-Host "Done updating table cells."